$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 196 (pushes existing rows 196-247 down to 197-248)
$ws.Rows.Item(196).Insert()

# Populate the newly inserted row 196 with the new weekly price-record data
$ws.Cells.Item(196, 1).Value = 1
$ws.Cells.Item(196, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(196, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(196, 4).Value = 44722
$ws.Cells.Item(196, 5).Value = 15
$ws.Cells.Item(196, 6).Value = "Fruta"
$ws.Cells.Item(196, 7).Value = 100102
$ws.Cells.Item(196, 8).Value = "Cítricos"
$ws.Cells.Item(196, 9).Value = 100102003
$ws.Cells.Item(196, 10).Value = "Limón"
$ws.Cells.Item(196, 11).Value = "Tahití"
$ws.Cells.Item(196, 12).Value = "Primera"
$ws.Cells.Item(196, 13).Value = 300
$ws.Cells.Item(196, 14).Value = 31000
$ws.Cells.Item(196, 15).Value = 32000
$ws.Cells.Item(196, 16).Value = 31500
$ws.Cells.Item(196, 17).Value = "$/caja 24 kilos"
$ws.Cells.Item(196, 18).Value = "Perú"
$ws.Cells.Item(196, 19).Value = 1312
$ws.Cells.Item(196, 20).Value = 24
